$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A1 with the new shared string value
$ws.Range("A1").Value = "gsdHGoi;xh"

# Update the active selection to A11 (matches saved view state)
$ws.Range("A11").Select()
